# Chapter 16 "계산 셰이더 (Compute Shader)" slide - clarify the pipeline-inclusion
# bullet with an "(X)" aside and fix the "자운"->"자원" typo in the GPU-resource
# bullet, splitting runs the same way PowerPoint does when you type in the
# middle of existing text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 3: "렌더링 파이프라인에 직접 포함X 파이프라인 옆에 따로 존재" ---
$para3 = $tr.Paragraphs(3, 1)
$base3 = $para3.Start

# "파이프라인에 직접 포함" (base3+4 .. +15) -> "파이프라인에 직접 " + "포함"
$tr.Characters($base3 + 4, 10).Text = "파이프라인에 직접 "
$tr.Characters($base3 + 14, 2).Text = "포함"

# Insert "(" right before "X " and give it its own run, then turn "X " into "X) "
$xRange = $tr.Characters($base3 + 16, 2)
[void]$xRange.InsertBefore("(")
$tr.Characters($base3 + 16, 1).Text = "("
$tr.Characters($base3 + 17, 2).Text = "X) "

# --- Paragraph 5: "GPU자원의 자료를 직접 읽거나 GPU 자운에 직접 자료를 기록 가능" ---
$para5 = $tr.Paragraphs(5, 1)
$base5 = $para5.Start

# "자운에 직접 자료를 기록 가능" -> "자원에 " + "직접 자료를 기록 가능" (also fixes the 자운/자원 typo)
$tr.Characters($base5 + 22, 4).Text = "자원에 "
